$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add header for new column H (copy formatting from G1's header style)
$ws.Range("H1").Value = "Justifications (if any)"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) | Out-Null

# Fill H2:H4 with "-" matching existing G2:G4 row style
$ws.Range("H2").Value = "-"
$ws.Range("H3").Value = "-"
$ws.Range("H4").Value = "-"
$ws.Range("G2:G4").Copy()
$ws.Range("H2:H4").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# Update selection to match diff (H1:H4 active cell H1)
$ws.Range("H1:H4").Select() | Out-Null

# Adjust column H width (engine quantizes to 1/6 steps; 23.71 is the closest
# achievable input that rounds to the target stored width of 24.5)
$ws.Range("H1").EntireColumn.ColumnWidth = 23.71
